$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 655
$ws1.Range("F4").Value = 261
$ws1.Range("F6").Value = 10098
$ws1.Range("F7").Value = 331
$ws1.Range("F8").Value = 915
$ws1.Range("F9").Value = 1257
$ws1.Range("F10").Value = 6070
$ws1.Range("F11").Value = 12
$ws1.Range("F12").Value = 416
$ws1.Range("F15").Value = 3104
$ws1.Range("F16").Value = 32
$ws1.Range("F18").Value = 601
$ws1.Range("F23").Value = 1542

# Sheet "全部类型" (sheet4): update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 655
$ws4.Range("F5").Value = 261
$ws4.Range("F7").Value = 10098
$ws4.Range("F8").Value = 331
$ws4.Range("F9").Value = 915
$ws4.Range("F10").Value = 1257
$ws4.Range("F11").Value = 6070
$ws4.Range("F12").Value = 12
$ws4.Range("F13").Value = 416
$ws4.Range("F16").Value = 3104
$ws4.Range("F17").Value = 32
$ws4.Range("F19").Value = 601
$ws4.Range("F24").Value = 1542
